$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value (serial 45205 = 2023-10-06) for every
# data row (rows 2 through 469). The update bumps that date forward by one day
# (serial 45206 = 2023-10-07) for all of them.
$firstRow = 2
$lastRow = 469

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45205) {
        $cell.Value2 = 45206
    }
}
